$wb = $excel.ActiveWorkbook
$wsMeta  = $wb.Worksheets.Item("Metadata")
$wsRules = $wb.Worksheets.Item("Rules")

# ---------------------------------------------------------------------------
# 1. Insert a new "Condition_3" column into the Rules table, between
#    Condition_2 and Action_1 (i.e. new worksheet column D), shifting the
#    former Action_1 / Action_2 columns one place to the right.
# ---------------------------------------------------------------------------
$wsRules.Columns.Item(4).Insert()

# Populate the new column's data (header, metadata rows, and the two rule rows)
$wsRules.Range("D1").Value = "Condition_3"
$wsRules.Range("D2").Value = "CabinClass"
$wsRules.Range("D3").Value = "Equals"
$wsRules.Range("D4").Value = "B"
$wsRules.Range("D5").Value = "E"

# Give the freshly inserted column a sensible width (close match of the
# original author's auto-fit result).
$wsRules.Columns.Item(4).ColumnWidth = 12.333333333333334

# Resize the "Rules" ListObject/table so it covers the new column.
$lo = $wsRules.ListObjects.Item(1)
$lo.Resize($wsRules.Range("A1:F15"))

# The resize step can mis-name columns, so make sure every header (and thus
# every table column name) is exactly right afterwards.
$wsRules.Range("A1").Value = "Index"
$wsRules.Range("B1").Value = "Condition_1"
$wsRules.Range("C1").Value = "Condition_2"
$wsRules.Range("D1").Value = "Condition_3"
$wsRules.Range("E1").Value = "Action_1"
$wsRules.Range("F1").Value = "Action_2"

# ---------------------------------------------------------------------------
# 2. Update the rule Index value in row 4 (7 -> 1).
# ---------------------------------------------------------------------------
$wsRules.Range("A4").Value = 1

# ---------------------------------------------------------------------------
# 3. Move the active sheet / selection: Metadata becomes the active tab
#    (selection E2), Rules keeps selection I5 but is no longer the active
#    tab.
# ---------------------------------------------------------------------------
$wsRules.Range("I5").Select()
$wsMeta.Activate()
$wsMeta.Range("E2").Select()
